# edit.ps1
# Applies the OOXML-level change described by the diff:
#  - Adds <w:ilvl w:val="0"/> into the <w:numPr> of the plain-text
#    (numId=0) paragraphs that were missing it.
#  - Removes the _GoBack bookmark from the end of the "B frame" bullet
#    paragraph and re-homes it (alone, in its own numId=0 paragraph)
#    after two newly authored bullet items are inserted:
#       "MP3" + "文件中的图片如专辑图片的编码格式是什么，现在还未知。"
#  - Adds the missing run of text to the final (numId=1) bullet
#    paragraph: "产生的时间戳是否要进行控制，如果是文件流，速度会很快，
#    现在录制的分割逻辑依据的是复用出来的packet的时间戳。"
#
# Strategy: Word's Range.InsertXML(xml) REPLACES the exact contents of
# the target Range with the supplied WordprocessingML, which lets us
# control w:numPr/w:ilvl, bookmarks and runs precisely -- something
# that plain text Find/Replace cannot do since it only ever touches
# the rendered text, never the paragraph properties.

function Wrap-Xml($bodyInner) {
    $pre = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $post = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pre + $bodyInner + $post
}

$d = $word.ActiveDocument

# --- 1) paragraphs 2-5 ("m3u8...", "HLS写视频...", "这样内存...", blank)
#        gain <w:ilvl w:val="0"/> ----------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p5 = $d.Paragraphs.Item(5)
$rngA = $d.Range($p2.Range.Start, $p5.Range.End)
$blockA = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>m3u8流，音频AAC去掉ADTS后少7个字节。</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>HLS写视频数据后多出来32个字节。</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>这样内存有泄漏。</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>'
$rngA.InsertXML((Wrap-Xml $blockA))

# --- 2) paragraphs 7-8 ("HLS流音频...", blank) gain <w:ilvl w:val="0"/> --
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$rngB = $d.Range($p7.Range.Start, $p8.Range.End)
$blockB = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>HLS流音频(AAC)数据送到FLVMuxer中难以识别。</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>'
$rngB.InsertXML((Wrap-Xml $blockB))

# --- 3) paragraphs 9-11 (the "B frame" bullet through the trailing
#        blank bullet) are replaced by five paragraphs: the "B frame"
#        bullet (bookmark removed), a blank ilvl=0 paragraph, the new
#        "MP3..." bullet, a new numId=0 paragraph holding the relocated
#        _GoBack bookmark, and the final bullet now carrying its text --
$p9 = $d.Paragraphs.Item(9)
$p11 = $d.Paragraphs.Item(11)
$rngC = $d.Range($p9.Range.Start, $p11.Range.End)
$blockC = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>拉有B帧的流时，时间戳如何处理，以为有可能DTS &lt; PTS，如果PTS=0，那么DTS就有可能是负数，这个时候ImportAVPacket方法中，Uint64改为int64是最好的解决方式。</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>MP3</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>文件中的图片如专辑图片的编码格式是什么，现在还未知。</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>产生的时间戳是否要进行控制，如果是文件流，速度会很快，现在录制的分割逻辑依据的是复用出来的packet的时间戳。</w:t></w:r></w:p>'
$rngC.InsertXML((Wrap-Xml $blockC))

Write-Host "Edit applied."
